# "Add files via upload" - appends one new menu item (vegetable-mayonnaise-sandwich)
# as row 32 of Sheet1, right after the existing "Peri-Peri French Frise" row (row 31).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate() | Out-Null

# New row: A = item name, C = price, D = image file name (column B stays blank,
# same layout used by every other row in the sheet).
$ws.Range("A32").Value = "vegetable-mayonnaise-sandwich"
$ws.Range("C32").Value = 20
$ws.Range("D32").Value = "vegetable-mayonnaise-sandwich.jpg"

# Copy the currency number format used by the rest of column C (style index 2)
# onto the new price cell so it renders/serializes the same way.
$ws.Range("C32").NumberFormat = $ws.Range("C31").NumberFormat

# Reflect the view state the workbook was left in after the edit: scrolled so
# row 16 is the first visible row, with C32 as the active/selected cell.
$excel.ActiveWindow.ScrollRow = 16
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("C32").Select() | Out-Null
